$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Remove the hyperlink that was previously on B1
$ws.Hyperlinks.Delete()

# Clear any special formatting (e.g. hyperlink style) on B1, reset to normal
$ws.Range("B1").Style = "Normal"

# Add second row of data
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = 123

# Update the active selection to match target (F12)
$ws.Range("F12").Select()
